# Rename "Product" sheet to "Products"
$wb = $excel.ActiveWorkbook
$wsProducts = $wb.Worksheets.Item("Product")
$wsProducts.Name = "Products"

# Update the header cells on the Products sheet (maSP/tenSP -> idProduct/nameProduct)
$wsProducts.Range("A1").Value = "idProduct"
$wsProducts.Range("B1").Value = "nameProduct"

# Autofit column B on the Products sheet (bestFit column width) and select a cell there
$wsProducts.Columns.Item(2).AutoFit() | Out-Null

# Update selection on the Users sheet (no longer the active tab)
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("I13").Select() | Out-Null

# Update selection on the Category sheet
$wsCategory = $wb.Worksheets.Item("Category")
$wsCategory.Range("E17").Select() | Out-Null

# Make Products the active sheet/tab and set its selection
$wsProducts.Activate() | Out-Null
$wsProducts.Range("L10").Select() | Out-Null
